{"js": "// Update the header date and the 25 two-digit multiplication problems\n// in the table, per the commit's regenerated problem set.\nconst replacements = [\n  [\"2025-12-02 Tuesday\", \"2025-12-03 Wednesday\"],\n  [\"17\u00d753=\", \"59\u00d745=\"],\n  [\"77\u00d774=\", \"51\u00d738=\"],\n  [\"70\u00d770=\", \"85\u00d736=\"],\n  [\"59\u00d764=\", \"91\u00d749=\"],\n  [\"84\u00d715=\", \"52\u00d797=\"],\n  [\"49\u00d786=\", \"20\u00d722=\"],\n  [\"95\u00d725=\", \"72\u00d713=\"],\n  [\"18\u00d750=\", \"64\u00d746=\"],\n  [\"64\u00d751=\", \"85\u00d724=\"],\n  [\"67\u00d715=\", \"61\u00d798=\"],\n  [\"64\u00d741=\", \"20\u00d723=\"],\n  [\"15\u00d740=\", \"28\u00d793=\"],\n  [\"95\u00d774=\", \"16\u00d786=\"],\n  [\"31\u00d737=\", \"94\u00d785=\"],\n  [\"58\u00d718=\", \"27\u00d769=\"],\n  [\"36\u00d718=\", \"81\u00d775=\"],\n  [\"42\u00d795=\", \"71\u00d755=\"],\n  [\"17\u00d725=\", \"81\u00d773=\"],\n  [\"49\u00d758=\", \"27\u00d719=\"],\n  [\"60\u00d735=\", \"33\u00d793=\"],\n  [\"64\u00d728=\", \"13\u00d754=\"],\n  [\"74\u00d714=\", \"17\u00d735=\"],\n  [\"21\u00d788=\", \"34\u00d728=\"],\n  [\"80\u00d727=\", \"49\u00d739=\"],\n  [\"80\u00d790=\", \"15\u00d745=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and the 25 two-digit multiplication problems\n# in the table, per the commit's regenerated problem set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-02 Tuesday\", \"2025-12-03 Wednesday\"),\n    @(\"17\u00d753=\", \"59\u00d745=\"),\n    @(\"77\u00d774=\", \"51\u00d738=\"),\n    @(\"70\u00d770=\", \"85\u00d736=\"),\n    @(\"59\u00d764=\", \"91\u00d749=\"),\n    @(\"84\u00d715=\", \"52\u00d797=\"),\n    @(\"49\u00d786=\", \"20\u00d722=\"),\n    @(\"95\u00d725=\", \"72\u00d713=\"),\n    @(\"18\u00d750=\", \"64\u00d746=\"),\n    @(\"64\u00d751=\", \"85\u00d724=\"),\n    @(\"67\u00d715=\", \"61\u00d798=\"),\n    @(\"64\u00d741=\", \"20\u00d723=\"),\n    @(\"15\u00d740=\", \"28\u00d793=\"),\n    @(\"95\u00d774=\", \"16\u00d786=\"),\n    @(\"31\u00d737=\", \"94\u00d785=\"),\n    @(\"58\u00d718=\", \"27\u00d769=\"),\n    @(\"36\u00d718=\", \"81\u00d775=\"),\n    @(\"42\u00d795=\", \"71\u00d755=\"),\n    @(\"17\u00d725=\", \"81\u00d773=\"),\n    @(\"49\u00d758=\", \"27\u00d719=\"),\n    @(\"60\u00d735=\", \"33\u00d793=\"),\n    @(\"64\u00d728=\", \"13\u00d754=\"),\n    @(\"74\u00d714=\", \"17\u00d735=\"),\n    @(\"21\u00d788=\", \"34\u00d728=\"),\n    @(\"80\u00d727=\", \"49\u00d739=\"),\n    @(\"80\u00d790=\", \"15\u00d745=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
